# Rename MinimalInformationGain to ItemSeverityBoundary
#
# The "Options" sheet (the active/second tab) holds a small lookup table in
# A1:B3. Cell A3 contained the label "Minimal information gain"; rename it
# to "Item severity boundary". Also bring over the cosmetic view changes
# that Excel recorded for this sheet after the edit: the selection moved to
# F6, and column A (bestFit) grew slightly to accommodate the new label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Options")

# Rename the lookup-table label.
$ws.Range("A3").Value = "Item severity boundary"

# Column A auto-fit to the new (slightly wider rendered) label.
$ws.Columns.Item(1).ColumnWidth = 21

# Selection left on the sheet after the edit.
$ws.Range("F6").Select() | Out-Null
